# Regenerate the handback status report: refresh the timestamp columns
# to reflect the latest run ("Generate Report for Handback").

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for fc9c9811-...md
# (this value is also mirrored on the de-de sheet's "Correspond Handoff
# Datetime" column, since it is the same underlying xliff-generation time)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-19 13:04:58"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-19 13:04:54"
$wsZhCn.Range("K2").Value = "2016-08-19 13:05:18"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-19 13:04:58"
$wsDeDe.Range("K2").Value = "2016-08-19 13:05:25"
